# Apply the commit's change: add a check for file existence before calling
# IsFileOpen in the therapist forms. This results in the 3PFormSheet
# (Rooms2 / "3PFormSheet") being reset/cleared and made the active sheet,
# and the "Last_Row_Read_3P" counter (Initials!D4) being bumped to 1.

$wb = $excel.ActiveWorkbook

# --- Update Initials!D4 (Last_Row_Read_3P) from 0 to 1 ---
$wsInitials = $wb.Worksheets.Item("Initials")
$wsInitials.Range("D4").Value = 1

# --- Clear the stale test data from 3PFormSheet rows 2-6, columns A:E ---
$ws3P = $wb.Worksheets.Item("3PFormSheet")
$ws3P.Range("A2:E6").ClearContents()

# --- Make 3PFormSheet the active / selected sheet (was 8PFormSheet) ---
$ws3P.Select()
$ws3P.Range("A2:E50").Select()

$wb.Save()
